# References/RAM.xlsx - "object per scanline counter"
# A new HRAM variable (FFB0 / curr_enemy_count) was inserted into the
# address table between FFAF and FFF0, pushing the existing FFF0..FFFE
# block down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row, shifting only the HRAM table (columns J:M) down by one
# -- the WRAM table on the left (columns A:G) is untouched -- so rows
# 45:54 move to 46:55 within J:M only.
$ws.Range("J45:M45").Insert(-4121)

# The shift-down nudges every row>=45 reference in the sheet's formulas,
# including the unrelated WRAM byte-count formula in F2; put it back.
$ws.Range("F2").Formula = "=SUM(D7:D999)"

# New row 45 becomes the new open/terminal slot (FFB0, no length, "(empty)").
$ws.Range("J45").Value = "FFB0"
$ws.Range("K45").Value = "FFEF"
$ws.Range("M45").Value = "(empty)"

# Row 44 used to be the open/terminal slot (FFAF, no length, "(empty)").
# It now becomes a real entry: 1 byte for curr_enemy_count.
$ws.Range("L44").Value = 1
$ws.Range("M44").Value = "curr_enemy_count"

# Re-enter the address formula across K7:K44 so the newly-filled K44 joins
# the existing fill as one contiguous shared formula.
$ws.Range("K7:K44").Formula = "=DEC2HEX(HEX2DEC(J7)+L7-1)"

# Re-enter the address formula across the remaining block so it forms its
# own shared-formula group.
$ws.Range("K46:K54").Formula = "=DEC2HEX(HEX2DEC(J46)+L46-1)"

# The scanned byte range grew by one row (L7:L999 -> L7:L1000).
$ws.Range("G2").Formula = "=SUM(L7:L1000)"

# Reflect the user's final selection/scroll position.
$ws.Range("G43").Select()
